# "Tried to implement Penality Reward System (unfinished)"
#
# 1. Weekly Quantity sheet: remove the first weekly data point
#    (Order Week 45116.99999999999 / Requested quantity 60), shifting the
#    remaining weeks up by one row (A1:B19 -> A1:B18).
# 2. Monthly Trend sheet: adjust the first month's requested quantity
#    from 130 down to 70 (penalty applied).

$wb = $excel.ActiveWorkbook

$weekly = $wb.Worksheets.Item("Weekly Quantity")
$weekly.Rows(2).Delete()

$monthly = $wb.Worksheets.Item("Monthly Trend")
$monthly.Range("B2").Value = 70
